$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.209.63'
$ws.Range('E2').Value = '  +1.18%  '
$ws.Range('D3').Value = '1.783.83'
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '226.45'
$ws.Range('E5').Value = '  +1.04%  '
$ws.Range('E6').Value = '  +0.36%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '31.95'
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('E9').Value = '  +1.20%  '
$ws.Range('E10').Value = '  +2.31%  '
$ws.Range('D12').Value = '2.040.15'
$ws.Range('E12').Value = '  +0.20%  '
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.06'
$ws.Range('E13').Value = '  -1.23%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.781.51'
$ws.Range('E14').Value = '  +0.37%  '
$ws.Range('D15').Value = '34.180.38'
$ws.Range('E15').Value = '  +1.01%  '
$ws.Range('E16').Value = '  +2.25%  '
$ws.Range('E17').Value = '  +1.55%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '67.90'
$ws.Range('E18').Value = '  +2.02%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '247.17'
$ws.Range('E19').Value = '  +3.75%  '
$ws.Range('D20').Value = '0.0₃0799'
$ws.Range('E20').Value = '  +3.47%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.99'
$ws.Range('E21').Value = '  +4.12%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.10'
$ws.Range('E23').Value = '  +2.44%  '
$ws.Range('E24').Value = '  -0.73%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '162.27'
$ws.Range('E25').Value = '  +1.36%  '
$ws.Range('E26').Value = '  +2.68%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.32'
$ws.Range('E27').Value = '  +1.60%  '
$ws.Range('E28').Value = '  +1.77%  '
$ws.Range('E29').Value = '  +0.18%  '
$ws.Range('E30').Value = '  +0.82%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0521'
$ws.Range('E31').Value = '  +2.15%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.75'
$ws.Range('E32').Value = '  +4.81%  '
$ws.Range('E33').Value = '  +5.72%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.80'
$ws.Range('E34').Value = '  -0.74%  '
$ws.Range('D35').Value = '1.444.78'
$ws.Range('E35').Value = '  +4.55%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.653'
$ws.Range('E36').Value = '  +2.84%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.41'
$ws.Range('E37').Value = '  +7.93%  '
$ws.Range('E38').Value = '  +3.91%  '
$ws.Range('E39').Value = '  +1.04%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.38'
$ws.Range('E40').Value = '  -0.61%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '80.35'
$ws.Range('E41').Value = '  +2.53%  '
$ws.Range('E42').Value = '  +1.81%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '13.67'
$ws.Range('E43').Value = '  +0.86%  '
$ws.Range('E44').Value = '  +0.80%  '
$ws.Range('E45').Value = '  +4.24%  '
$ws.Range('E46').Value = '  +0.50%  '
$ws.Range('E47').Value = '  -0.21%  '
$ws.Range('E48').Value = '  -4.94%  '
$ws.Range('D49').Value = '1.942.49'
$ws.Range('E49').Value = '  +0.25%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '104.73'
$ws.Range('E50').Value = '  -2.67%  '
$ws.Range('E51').Value = '  +0.12%  '
